$d = $word.ActiveDocument

# 1) Replace the tail of the "most difficult part of the project" sentence:
#    "...Graphical User Interface (GUI) as well as the web scraping component. "
#    -> "...Graphical User Interface (GUI) where QT Creator was used."
$d.Content.Find.Execute(
    "as well as the web scraping component. ", $false, $false, $false, $false,
    $false, $true, 1, $false, "where QT Creator was used.", 2) | Out-Null

# 2) Rewrite the "web scraping component" sentence to talk about restaurants instead:
#    "The web scraping component of this project relies on a website"
#    -> "The restaurants that this project relies on is a website"
$d.Content.Find.Execute(
    "The web scraping component of this project relies on a website",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The restaurants that this project relies on is a website", 2) | Out-Null

# 3) Fix subject/verb agreement: "drastically is small" -> "drastically are small"
$d.Content.Find.Execute(
    "on website changes drastically is small", $false, $false, $false, $false,
    $false, $true, 1, $false, "on website changes drastically are small", 2) | Out-Null
